{"js": "// Replacement list: [oldText, newText] pairs, in document (reading) order \u2014\n// one entry per table cell, matching the order the cells appear in the body.\nconst replacements = [[\"42+32=\", \"87-22=\"], [\"14+17=\", \"79+6=\"], [\"57-54=\", \"39+36=\"], [\"86-14=\", \"60-8=\"], [\"50-46=\", \"59-20=\"], [\"23+32=\", \"92-17=\"], [\"31+7=\", \"78+14=\"], [\"66-53=\", \"15-3=\"], [\"47-9=\", \"27+7=\"], [\"94-11=\", \"58-43=\"], [\"87-72=\", \"7+5=\"], [\"76-23=\", \"96-2=\"], [\"4+19=\", \"52-10=\"], [\"89-24=\", \"68-31=\"], [\"13+85=\", \"70-67=\"], [\"54-13=\", \"19-10=\"], [\"47+8=\", \"40+40=\"], [\"46+34=\", \"44+18=\"], [\"64-33=\", \"91-67=\"], [\"76-5=\", \"43-42=\"], [\"23+74=\", \"35-7=\"], [\"81-27=\", \"29-16=\"], [\"81-8=\", \"69-52=\"], [\"30+42=\", \"39+26=\"], [\"50-11=\", \"33-21=\"], [\"50-15=\", \"76-50=\"], [\"38+38=\", \"64-3=\"], [\"96-3=\", \"19+2=\"], [\"53-1=\", \"65+5=\"], [\"20+53=\", \"22+32=\"], [\"81-62=\", \"48+6=\"], [\"95-59=\", \"50-13=\"], [\"66-3=\", \"96-68=\"], [\"82-25=\", \"11+23=\"], [\"81+2=\", \"8+51=\"], [\"82-82=\", \"27+70=\"], [\"12+57=\", \"45+22=\"], [\"42-35=\", \"64+3=\"], [\"82+15=\", \"35-15=\"], [\"26-0=\", \"98-82=\"], [\"50-49=\", \"27+35=\"], [\"44+46=\", \"23+62=\"], [\"24+22=\", \"44-29=\"], [\"52-2=\", \"31+64=\"], [\"22+2=\", \"19+46=\"], [\"11+72=\", \"55+39=\"], [\"82-35=\", \"4+88=\"], [\"50+35=\", \"35-27=\"], [\"4+75=\", \"80-17=\"], [\"63-60=\", \"50-20=\"], [\"10+55=\", \"36-32=\"], [\"62-2=\", \"47+23=\"], [\"30-1=\", \"15+70=\"], [\"40+53=\", \"20+54=\"], [\"53+23=\", \"17-3=\"], [\"83-66=\", \"53+6=\"], [\"2+43=\", \"48-20=\"], [\"77+0=\", \"61-4=\"], [\"51+40=\", \"44+51=\"], [\"64+9=\", \"98-0=\"], [\"16+82=\", \"4+68=\"], [\"89-64=\", \"2+75=\"], [\"53-45=\", \"11-4=\"], [\"83-47=\", \"70+17=\"], [\"84-30=\", \"95-1=\"], [\"52-21=\", \"61-0=\"], [\"57-23=\", \"42+48=\"], [\"68-54=\", \"92-3=\"], [\"37+31=\", \"15-11=\"], [\"70-15=\", \"95+3=\"], [\"66+20=\", \"46-20=\"], [\"19+41=\", \"49-22=\"], [\"5+68=\", \"89-7=\"], [\"13+51=\", \"79-25=\"], [\"78-23=\", \"42+43=\"], [\"34+58=\", \"75-30=\"], [\"23+12=\", \"28+54=\"], [\"37-1=\", \"48-43=\"], [\"22+39=\", \"20+40=\"], [\"10+15=\", \"86-67=\"], [\"77-1=\", \"79-67=\"], [\"36+5=\", \"86-48=\"], [\"81-76=\", \"24+26=\"], [\"5+61=\", \"23+71=\"], [\"10+87=\", \"18+26=\"], [\"38-20=\", \"21+10=\"], [\"19+9=\", \"28+49=\"], [\"24-13=\", \"76+19=\"], [\"18-12=\", \"11+59=\"], [\"34-0=\", \"19+46=\"], [\"57+40=\", \"93-34=\"], [\"71-68=\", \"83-38=\"], [\"8+48=\", \"31-18=\"], [\"40+27=\", \"82-1=\"], [\"80+9=\", \"73-69=\"], [\"84-76=\", \"63-47=\"], [\"32+8=\", \"7+84=\"], [\"94-76=\", \"30-11=\"], [\"77-16=\", \"3+36=\"], [\"43+37=\", \"11+23=\"]];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values[0].length;\n\n// Primary pass: walk the table in reading order (row-major) and match each\n// cell against the replacement list at the same position.\nlet i = 0;\nconst unresolved = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    if (i >= replacements.length) break;\n    const [oldText, newText] = replacements[i];\n    const cell = table.getCell(r, c);\n    cell.load(\"value\");\n    await context.sync();\n    if (cell.value === oldText) {\n      cell.value = newText;\n    } else {\n      unresolved.push([oldText, newText]);\n    }\n    i++;\n  }\n}\nawait context.sync();\n\n// Fallback pass: if a cell's current text didn't line up positionally with\n// the expected old value (e.g. table layout differs from assumption),\n// search every cell for an exact-text match and fix it there instead.\nif (unresolved.length > 0) {\n  for (const [oldText, newText] of unresolved) {\n    outer:\n    for (let r = 0; r < table.rowCount; r++) {\n      for (let c = 0; c < columnCount; c++) {\n        const cell = table.getCell(r, c);\n        cell.load(\"value\");\n        await context.sync();\n        if (cell.value === oldText) {\n          cell.value = newText;\n          break outer;\n        }\n      }\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Replacement list: (oldText, newText) pairs, in document (reading) order \u2014\n# one entry per table cell, matching the order the cells appear in the table.\n$replacements = @(\n    ,@(\"42+32=\", \"87-22=\")\n    ,@(\"14+17=\", \"79+6=\")\n    ,@(\"57-54=\", \"39+36=\")\n    ,@(\"86-14=\", \"60-8=\")\n    ,@(\"50-46=\", \"59-20=\")\n    ,@(\"23+32=\", \"92-17=\")\n    ,@(\"31+7=\", \"78+14=\")\n    ,@(\"66-53=\", \"15-3=\")\n    ,@(\"47-9=\", \"27+7=\")\n    ,@(\"94-11=\", \"58-43=\")\n    ,@(\"87-72=\", \"7+5=\")\n    ,@(\"76-23=\", \"96-2=\")\n    ,@(\"4+19=\", \"52-10=\")\n    ,@(\"89-24=\", \"68-31=\")\n    ,@(\"13+85=\", \"70-67=\")\n    ,@(\"54-13=\", \"19-10=\")\n    ,@(\"47+8=\", \"40+40=\")\n    ,@(\"46+34=\", \"44+18=\")\n    ,@(\"64-33=\", \"91-67=\")\n    ,@(\"76-5=\", \"43-42=\")\n    ,@(\"23+74=\", \"35-7=\")\n    ,@(\"81-27=\", \"29-16=\")\n    ,@(\"81-8=\", \"69-52=\")\n    ,@(\"30+42=\", \"39+26=\")\n    ,@(\"50-11=\", \"33-21=\")\n    ,@(\"50-15=\", \"76-50=\")\n    ,@(\"38+38=\", \"64-3=\")\n    ,@(\"96-3=\", \"19+2=\")\n    ,@(\"53-1=\", \"65+5=\")\n    ,@(\"20+53=\", \"22+32=\")\n    ,@(\"81-62=\", \"48+6=\")\n    ,@(\"95-59=\", \"50-13=\")\n    ,@(\"66-3=\", \"96-68=\")\n    ,@(\"82-25=\", \"11+23=\")\n    ,@(\"81+2=\", \"8+51=\")\n    ,@(\"82-82=\", \"27+70=\")\n    ,@(\"12+57=\", \"45+22=\")\n    ,@(\"42-35=\", \"64+3=\")\n    ,@(\"82+15=\", \"35-15=\")\n    ,@(\"26-0=\", \"98-82=\")\n    ,@(\"50-49=\", \"27+35=\")\n    ,@(\"44+46=\", \"23+62=\")\n    ,@(\"24+22=\", \"44-29=\")\n    ,@(\"52-2=\", \"31+64=\")\n    ,@(\"22+2=\", \"19+46=\")\n    ,@(\"11+72=\", \"55+39=\")\n    ,@(\"82-35=\", \"4+88=\")\n    ,@(\"50+35=\", \"35-27=\")\n    ,@(\"4+75=\", \"80-17=\")\n    ,@(\"63-60=\", \"50-20=\")\n    ,@(\"10+55=\", \"36-32=\")\n    ,@(\"62-2=\", \"47+23=\")\n    ,@(\"30-1=\", \"15+70=\")\n    ,@(\"40+53=\", \"20+54=\")\n    ,@(\"53+23=\", \"17-3=\")\n    ,@(\"83-66=\", \"53+6=\")\n    ,@(\"2+43=\", \"48-20=\")\n    ,@(\"77+0=\", \"61-4=\")\n    ,@(\"51+40=\", \"44+51=\")\n    ,@(\"64+9=\", \"98-0=\")\n    ,@(\"16+82=\", \"4+68=\")\n    ,@(\"89-64=\", \"2+75=\")\n    ,@(\"53-45=\", \"11-4=\")\n    ,@(\"83-47=\", \"70+17=\")\n    ,@(\"84-30=\", \"95-1=\")\n    ,@(\"52-21=\", \"61-0=\")\n    ,@(\"57-23=\", \"42+48=\")\n    ,@(\"68-54=\", \"92-3=\")\n    ,@(\"37+31=\", \"15-11=\")\n    ,@(\"70-15=\", \"95+3=\")\n    ,@(\"66+20=\", \"46-20=\")\n    ,@(\"19+41=\", \"49-22=\")\n    ,@(\"5+68=\", \"89-7=\")\n    ,@(\"13+51=\", \"79-25=\")\n    ,@(\"78-23=\", \"42+43=\")\n    ,@(\"34+58=\", \"75-30=\")\n    ,@(\"23+12=\", \"28+54=\")\n    ,@(\"37-1=\", \"48-43=\")\n    ,@(\"22+39=\", \"20+40=\")\n    ,@(\"10+15=\", \"86-67=\")\n    ,@(\"77-1=\", \"79-67=\")\n    ,@(\"36+5=\", \"86-48=\")\n    ,@(\"81-76=\", \"24+26=\")\n    ,@(\"5+61=\", \"23+71=\")\n    ,@(\"10+87=\", \"18+26=\")\n    ,@(\"38-20=\", \"21+10=\")\n    ,@(\"19+9=\", \"28+49=\")\n    ,@(\"24-13=\", \"76+19=\")\n    ,@(\"18-12=\", \"11+59=\")\n    ,@(\"34-0=\", \"19+46=\")\n    ,@(\"57+40=\", \"93-34=\")\n    ,@(\"71-68=\", \"83-38=\")\n    ,@(\"8+48=\", \"31-18=\")\n    ,@(\"40+27=\", \"82-1=\")\n    ,@(\"80+9=\", \"73-69=\")\n    ,@(\"84-76=\", \"63-47=\")\n    ,@(\"32+8=\", \"7+84=\")\n    ,@(\"94-76=\", \"30-11=\")\n    ,@(\"77-16=\", \"3+36=\")\n    ,@(\"43+37=\", \"11+23=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n# Primary pass: walk the table in reading order (row-major) and match each\n# cell against the replacement list at the same position.\n$unresolved = New-Object System.Collections.ArrayList\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($i -ge $replacements.Count) { break }\n    $oldText = $replacements[$i][0]\n    $newText = $replacements[$i][1]\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    $cellRange.MoveEnd(1, -1) | Out-Null\n    if ($cellRange.Text -eq $oldText) {\n      $cellRange.Text = $newText\n    } else {\n      [void]$unresolved.Add(@($oldText, $newText))\n    }\n    $i++\n  }\n}\n\n# Fallback pass: if a cell's current text didn't line up positionally with\n# the expected old value (e.g. table layout differs from assumption),\n# search every cell for an exact-text match and fix it there instead.\nforeach ($pair in $unresolved) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $found = $false\n  for ($r = 1; $r -le $rowCount -and -not $found; $r++) {\n    for ($c = 1; $c -le $colCount -and -not $found; $c++) {\n      $cell = $t.Cell($r, $c)\n      $cellRange = $cell.Range\n      $cellRange.MoveEnd(1, -1) | Out-Null\n      if ($cellRange.Text -eq $oldText) {\n        $cellRange.Text = $newText\n        $found = $true\n      }\n    }\n  }\n}\n"}
